$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rev. D -> Rev. E
$ws.Name = "Domino Single USB Rev. E"

# Re-assert the print area under the new sheet name so the workbook-level
# _xlnm.Print_Area defined name follows the rename too.
$ws.PageSetup.PrintArea = "A1:I14"

# Column width tweaks (silkscreen font ratio -> 20%)
$ws.Columns.Item(1).ColumnWidth = 4.16666666666667
$ws.Columns.Item(2).ColumnWidth = 4.16666666666667
$ws.Columns.Item(3).ColumnWidth = 26
$ws.Columns.Item(4).ColumnWidth = 27.8333333333333
$ws.Columns.Item(5).ColumnWidth = 31
$ws.Columns.Item(6).ColumnWidth = 26
$ws.Columns.Item(7).ColumnWidth = 40.1666666666667
$ws.Columns.Item(8).ColumnWidth = 62.3333333333333
$ws.Columns.Item(9).ColumnWidth = 24.1666666666667

# Row 15 gets an explicit (custom) height, same value as before
$ws.Rows.Item(15).RowHeight = 12.1
